$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ghi cong")
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 0

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("H13").Select()
